# Generate Report for Handoff
# Rotates the localization-status workbook's tracked source file from the
# old "08a78c31-a338-481e-82ca-38495c11ef11.md" GUID-named markdown file to
# the newly generated "e4685e38-6d26-49b2-8134-870123067cee.md", and bumps
# the handoff/handback timestamps that were recorded for the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "08a78c31-a338-481e-82ca-38495c11ef11"
$newGuid = "e4685e38-6d26-49b2-8134-870123067cee"

# The hyperlinks on A2/B2 point at this (unchanged) GitHub blob URL - only
# the visible/display text is refreshed to the new file name.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9555d8c4ee67a4328d1a23a80e8e6da83cd770de/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 08:57:25"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.a5c02cb4448193b5af9731960742a3843cd50713.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-16 08:57:19"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.a5c02cb4448193b5af9731960742a3843cd50713.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-16 08:57:25"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
